$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 53 holds the new day's entry: 2026/01/02, 逃离鸭科夫, 1127.
# Column A stores the date as literal text (like all prior rows), not an
# auto-converted date serial. Writing a formula that returns the literal
# string, then pasting its computed value back in place, sidesteps Excel's
# "looks like a date" auto-detection without registering a throwaway
# number-format style in the workbook.
$ws.Cells.Item(53, 1).Formula = "=""2026/01/02"""
$ws.Cells.Item(53, 1).Copy()
$ws.Cells.Item(53, 1).PasteSpecial(-4163)  # xlPasteValues

$ws.Cells.Item(53, 2).Value = "逃离鸭科夫"
$ws.Cells.Item(53, 3).Value = 1127

# Match the existing data rows' look (centered alignment) by copying the
# format from the row above, same as the rest of the table (rows 3-52).
$ws.Range("A52:C52").Copy()
$ws.Range("A53:C53").PasteSpecial(-4122)  # xlPasteFormats

$ws.Application.CutCopyMode = 0
